# Add 2022-Q3 data
# 1) Insert a new worksheet "2022-Q3" right after the "总计" sheet, populated
#    with the fund-holding detail table for that quarter (mirrors the layout
#    of the other quarterly sheets).
# 2) Update the "总计" (summary) sheet: shift the existing quarter rows down
#    by one and insert a new first data row for "2022-Q3".

$wb = $excel.ActiveWorkbook
$zj = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$newWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $zj)
$newWs.Name = "2022-Q3"

$headerStyleSrc = $zj.Cells.Item(1, 2)   # bold/centered header style (s=2)
$indexStyleSrc  = $zj.Cells.Item(2, 1)   # bold/centered index-column style (s=2)

$headerCols = @(2, 3, 4, 5, 6, 7, 8)
foreach ($c in $headerCols) {
    $headerStyleSrc.Copy($newWs.Cells.Item(1, $c))
}

$newWs.Cells.Item(1, 2).Value = "基金代码"
$newWs.Cells.Item(1, 3).Value = "基金名称"
$newWs.Cells.Item(1, 4).Value = "基金规模"
$newWs.Cells.Item(1, 5).Value = "股票总仓位"
$newWs.Cells.Item(1, 6).Value = "仓位占比"
$newWs.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newWs.Cells.Item(1, 8).Value = "仓位排名"

for ($r = 2; $r -le 6; $r++) {
    $indexStyleSrc.Copy($newWs.Cells.Item($r, 1))
}

$newWs.Cells.Item(2, 1).Value = 0
$newWs.Cells.Item(2, 2).Value = "'006102"
$newWs.Cells.Item(2, 3).Value = "浙商丰利增强债券"
$newWs.Cells.Item(2, 4).Value = "'93.25"
$newWs.Cells.Item(2, 5).Value = "'44.13"
$newWs.Cells.Item(2, 6).Value = "'1.55"
$newWs.Cells.Item(2, 7).Value = "'1.4454"
$newWs.Cells.Item(2, 8).Value = 10

$newWs.Cells.Item(3, 1).Value = 1
$newWs.Cells.Item(3, 2).Value = "'011179"
$newWs.Cells.Item(3, 3).Value = "浙商智选食品饮料股票A"
$newWs.Cells.Item(3, 4).Value = "'0.14"
$newWs.Cells.Item(3, 5).Value = "'91.42"
$newWs.Cells.Item(3, 6).Value = "'6.49"
$newWs.Cells.Item(3, 7).Value = "'0.0091"
$newWs.Cells.Item(3, 8).Value = 7

$newWs.Cells.Item(4, 1).Value = 2
$newWs.Cells.Item(4, 2).Value = "'011180"
$newWs.Cells.Item(4, 3).Value = "浙商智选食品饮料股票C"
$newWs.Cells.Item(4, 4).Value = "'0.08"
$newWs.Cells.Item(4, 5).Value = "'91.42"
$newWs.Cells.Item(4, 6).Value = "'6.49"
$newWs.Cells.Item(4, 7).Value = "'0.0052"
$newWs.Cells.Item(4, 8).Value = 7

$newWs.Cells.Item(5, 1).Value = 3
$newWs.Cells.Item(5, 2).Value = "'005429"
$newWs.Cells.Item(5, 3).Value = "渤海汇金睿选混合A"
$newWs.Cells.Item(5, 4).Value = "'0.13"
$newWs.Cells.Item(5, 5).Value = "'30.88"
$newWs.Cells.Item(5, 6).Value = "'1.34"
$newWs.Cells.Item(5, 7).Value = "'0.0017"
$newWs.Cells.Item(5, 8).Value = 8

$newWs.Cells.Item(6, 1).Value = 4
$newWs.Cells.Item(6, 2).Value = "'005430"
$newWs.Cells.Item(6, 3).Value = "渤海汇金睿选混合C"
$newWs.Cells.Item(6, 4).Value = "'0.01"
$newWs.Cells.Item(6, 5).Value = "'30.88"
$newWs.Cells.Item(6, 6).Value = "'1.34"
$newWs.Cells.Item(6, 7).Value = "'0.0001"
$newWs.Cells.Item(6, 8).Value = 8

# ---------------------------------------------------------------------
# Step 2: update the "总计" (summary) sheet - shift rows down and insert
# the new 2022-Q3 summary row at the top of the data (row 2).
# ---------------------------------------------------------------------

# Row 9 ("2020-Q4") is brand-new; clone the index-column style from row 8
# before writing values into it.
$zj.Cells.Item(8, 1).Copy($zj.Cells.Item(9, 1))

$zj.Cells.Item(9, 1).Value = 7
$zj.Cells.Item(9, 2).Value = "2020-Q4"
$zj.Cells.Item(9, 3).Value = 7
$zj.Cells.Item(9, 4).Value = 0.12

$zj.Cells.Item(8, 1).Value = 6
$zj.Cells.Item(8, 2).Value = "2021-Q1"
$zj.Cells.Item(8, 3).Value = 10
$zj.Cells.Item(8, 4).Value = 1.08

$zj.Cells.Item(7, 1).Value = 5
$zj.Cells.Item(7, 2).Value = "2021-Q2"
$zj.Cells.Item(7, 3).Value = 6
$zj.Cells.Item(7, 4).Value = 1.11

$zj.Cells.Item(6, 1).Value = 4
$zj.Cells.Item(6, 2).Value = "2021-Q3"
$zj.Cells.Item(6, 3).Value = 20
$zj.Cells.Item(6, 4).Value = 4.26

$zj.Cells.Item(5, 1).Value = 3
$zj.Cells.Item(5, 2).Value = "2021-Q4"
$zj.Cells.Item(5, 3).Value = 7
$zj.Cells.Item(5, 4).Value = 1.7

$zj.Cells.Item(4, 1).Value = 2
$zj.Cells.Item(4, 2).Value = "2022-Q1"
$zj.Cells.Item(4, 3).Value = 18
$zj.Cells.Item(4, 4).Value = 1.83

$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(3, 2).Value = "2022-Q2"
$zj.Cells.Item(3, 3).Value = 8
$zj.Cells.Item(3, 4).Value = 1.68

$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 2).Value = "2022-Q3"
$zj.Cells.Item(2, 3).Value = 5
$zj.Cells.Item(2, 4).Value = 1.46
